# BrandDetails.xlsx update:
#  - keep only the header row + first data row (rows 3-5 removed)
#  - update the mobile number on the remaining data row
#  - update the JSON "Opertor" payload to match the new operator/number/id
#  - widen column I slightly
#  - refresh the selection to match the shrunk used range

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three trailing data rows (rows 3,4,5), leaving header (row1) + row2.
$ws.Rows("3:5").Delete() | Out-Null

# Row 2: new mobile number.
$ws.Range("E2").Value = 9652356895

# Row 2: "Opertor" column JSON payload now reflects Airtel / id 1 / new number.
$ws.Range("I2").Value = '[{"id":1,"operator":"Airtel","mobile_number":"9652356895","phonecode":91}]'

# Column I a touch wider (displayed width 85 -> 88).
$ws.Columns("I").ColumnWidth = 87.16666666666667

# Selection now only spans the remaining two rows.
$ws.Range("F1:F2").Select() | Out-Null
